$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated the raw ignition readings in column A (values re-measured/rounded).
# Column B keeps its existing "=A*2-20" formulas and just recalculates.
$newValues = @(2, 7, 9, 11, 12, 13, 14, 15, 16, 16, 17, 17, 18, 18, 19, 19, 20)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# Selection ends up on the recalculated detonation (B) column.
$ws.Range("B1:B17").Select()
